# Table 1 touch-up: the "Beginning Date" value in B5 was still stored as a
# literal date (1/1/1970), but the rest of the date column (B6, "End Date")
# is entered as free text so it can carry the custom "Month d, yyyy" label.
# Replace B5 with the matching text label (keeping the same typo style
# already present in B6, e.g. "Decmeber") so both date cells are consistent.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Janurary 1, 1970"

# Leave the selection where the author last left it before saving.
$null = $ws.Range("B2").Select()
